$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last data row (row 64), which no longer exists in the refreshed dataset
$ws.Rows(64).Delete()

# Re-generate the dataset (index/address/price/url) for rows 2-63 from the refreshed scrape
$addresses = @(
    '1085 103rd Ave NE, Bellevue, WA',
    'Park 433 | 433 Bellevue Way SE, Bellevue, WA',
    '537 Bellevue Way SE #1697806, Bellevue, WA 98004',
    'Alley 111 | 11011 NE 9th St, Bellevue, WA',
    '177 107th Ave NE, Bellevue, WA',
    'Aventine | 211 112th Ave NE, Bellevue, WA',
    '12 Central Square | 10290 NE 12th St, Bellevue, WA',
    'Cerasa | 10961 NE 2nd Pl, Bellevue, WA',
    'City Square Bellevue | 938 110th Ave NE, Bellevue, WA',
    '10245 Main St #73JGW1E7W, Bellevue, WA 98004',
    'Sylva on Main Apartments | 10701 Main St, Bellevue, WA',
    'Soma Towers | 288 106th Ave NE, Bellevue, WA',
    '88 102nd Ave NE #161KHGUNZ, Bellevue, WA 98004',
    'Avalon Towers Bellevue | 10349 NE 10th St, Bellevue, WA',
    'Avalon Bellevue | 11000 NE 10th St, Bellevue, WA',
    'AMLI Bellevue Park | 10001 NE 1st St, Bellevue, WA',
    'Metro 112 Apartments | 317 112th Ave NE, Bellevue, WA',
    'Main Street Flats | 10575 Main St, Bellevue, WA',
    'The Bravern | 688 110th Ave NE, Bellevue, WA',
    'Borgata | 37 103rd Ave NE, Bellevue, WA',
    'Ashton Bellevue | 10710 NE 10th St, Bellevue, WA',
    'Two Lincoln Tower | 10485 NE 6th St, Bellevue, WA',
    '1515 Bellevue Way NE #3375777, Bellevue, WA 98004',
    'Venn at Main | 10333 NE 1st St, Bellevue, WA',
    'The Meyden | 10333 Main St, Bellevue, WA',
    'Avalon Meydenbauer | 10410 NE 2nd St, Bellevue, WA',
    'Lux | 1000 100th Ave NE, Bellevue, WA',
    'Brio | 11130 NE 10th St, Bellevue, WA',
    'BLU | 75 102nd Ave NE, Bellevue, WA',
    'Elements Apartments | 958 111th Ave NE, Bellevue, WA',
    'TEN20 | 1020 108th Ave NE, Bellevue, WA',
    '11101 NE 12th St APT 107, Bellevue, WA 98004',
    '900 108th Ave NE | 900 108th Ave NE, Bellevue, WA',
    '300 110th Ave NE, Bellevue, WA',
    '10305 NE 16th St APT I7, Bellevue, WA 98004',
    '410 102nd Ave SE, Bellevue, WA 98004',
    '118 107th Ave NE #315, Bellevue, WA 98004',
    '9922 Lake Washington Blvd NE, Bellevue, WA',
    '650 Bellevue Way NE, Bellevue, WA',
    '10610 NE 9th Pl, Bellevue, WA',
    '205 105th Ave SE, Bellevue, WA',
    '10650 NE 9th Pl UNIT 1524, Bellevue, WA 98004',
    '1100 106th Ave NE APT 606, Bellevue, WA 98004',
    '1620 103rd Ave NE, Bellevue, WA 98004',
    '188 Bellevue Way NE, Bellevue, WA',
    '1188 106th Ave NE APT 321, Bellevue, WA 98004',
    '111 108th Ave NE UNIT A411, Bellevue, WA 98004',
    '812 100th Ave NE, Bellevue, WA',
    '10700 NE 4th St UNIT 428, Bellevue, WA 98004',
    '10008 NE 16th Pl #ADU, Bellevue, WA 98004',
    '424 102nd Ave SE APT 308, Bellevue, WA 98004',
    '550 100th Ave #16, Bellevue, WA 98004',
    '10226 SE 6th Street | 10226 SE 6th St, Bellevue, WA',
    '10042 Main St APT 402, Bellevue, WA 98004',
    '10608 NE 2nd St FLOOR 12-ID117, Bellevue, WA 98004',
    '10201 SE 3rd St #1, Bellevue, WA 98004',
    '425 Bellevue Way SE #25, Bellevue, WA 98004',
    '425 Bellevue Way SE #62, Bellevue, WA 98004',
    '417 99th Ave NE #G, Bellevue, WA 98004',
    '909 112th Ave NE, Bellevue, WA',
    '10608 NE 2nd St, Bellevue, WA',
    '125 108th Ave SE APT 12, Bellevue, WA 98004'
)

$prices = @(
    1795,
    1795,
    1650,
    1888,
    2450,
    1736,
    2671,
    1970,
    1670,
    2000,
    1627,
    1915,
    2250,
    1940,
    1790,
    2035,
    1875,
    1710,
    2165,
    1855,
    3299,
    3902,
    2455,
    2355,
    1788,
    1765,
    2850,
    2265,
    2595,
    1761,
    3105,
    1780,
    1980,
    2200,
    3000,
    2100,
    2795,
    1550,
    4800,
    2700,
    2220,
    3800,
    1800,
    2649,
    4500,
    1900,
    1750,
    1925,
    3100,
    2000,
    2700,
    1895,
    1495,
    3995,
    2830,
    1275,
    1400,
    1325,
    2200,
    3560,
    2830,
    1750
)

$urls = @(
    'https://www.zillow.com/b/1085-103rd-ave-ne-bellevue-wa-BG9pnD/',
    'https://www.zillow.com/b/park-433-bellevue-wa-9VSJsn/',
    'https://www.zillow.com/homedetails/537-Bellevue-Way-SE-1697806-Bellevue-WA-98004/2098531044_zpid/',
    'https://www.zillow.com/b/alley-111-bellevue-wa-5zpjzy/',
    'https://www.zillow.com/b/177-107th-ave-ne-bellevue-wa-5Xrs4s/',
    'https://www.zillow.com/b/aventine-bellevue-wa-5XjKwL/',
    'https://www.zillow.com/b/12-central-square-bellevue-wa-5XjFqC/',
    'https://www.zillow.com/b/cerasa-bellevue-wa-BKR5YY/',
    'https://www.zillow.com/b/city-square-bellevue-bellevue-wa-5XkJ93/',
    'https://www.zillow.com/homedetails/10245-Main-St-73JGW1E7W-Bellevue-WA-98004/2073066244_zpid/',
    'https://www.zillow.com/b/sylva-on-main-apartments-bellevue-wa-5Xnhwr/',
    'https://www.zillow.com/b/soma-towers-bellevue-wa-5bR2bQ/',
    'https://www.zillow.com/homedetails/88-102nd-Ave-NE-161KHGUNZ-Bellevue-WA-98004/2070452484_zpid/',
    'https://www.zillow.com/b/avalon-towers-bellevue-bellevue-wa-5XjFg5/',
    'https://www.zillow.com/b/avalon-bellevue-bellevue-wa-5XjPfD/',
    'https://www.zillow.com/b/amli-bellevue-park-bellevue-wa-5Xj9VV/',
    'https://www.zillow.com/b/metro-112-apartments-bellevue-wa-5Xk5rF/',
    'https://www.zillow.com/b/main-street-flats-bellevue-wa-65ZSwC/',
    'https://www.zillow.com/b/the-bravern-bellevue-wa-5XjRT2/',
    'https://www.zillow.com/b/borgata-bellevue-wa-5XjRc4/',
    'https://www.zillow.com/b/ashton-bellevue-bellevue-wa-5XjV9N/',
    'https://www.zillow.com/b/two-lincoln-tower-bellevue-wa-BHmS47/',
    'https://www.zillow.com/homedetails/1515-Bellevue-Way-NE-3375777-Bellevue-WA-98004/2079933369_zpid/',
    'https://www.zillow.com/b/venn-at-main-bellevue-wa-9kBw6T/',
    'https://www.zillow.com/b/the-meyden-bellevue-wa-5hJ4p9/',
    'https://www.zillow.com/b/avalon-meydenbauer-bellevue-wa-5XjGYw/',
    'https://www.zillow.com/b/lux-bellevue-wa-5Znksy/',
    'https://www.zillow.com/b/brio-bellevue-wa-BW9wC7/',
    'https://www.zillow.com/b/blu-bellevue-wa-BcxTVp/',
    'https://www.zillow.com/b/elements-apartments-bellevue-wa-5XjW4s/',
    'https://www.zillow.com/b/ten20-bellevue-wa-5XjQBF/',
    'https://www.zillow.com/homedetails/11101-NE-12th-St-APT-107-Bellevue-WA-98004/2082985886_zpid/',
    'https://www.zillow.com/b/900-108th-ave-ne-bellevue-wa-5YJVG4/',
    'https://www.zillow.com/b/300-110th-ave-ne-bellevue-wa-5XkXpv/',
    'https://www.zillow.com/homedetails/10305-NE-16th-St-APT-I7-Bellevue-WA-98004/2070238834_zpid/',
    'https://www.zillow.com/homedetails/410-102nd-Ave-SE-Bellevue-WA-98004/2075772896_zpid/',
    'https://www.zillow.com/homedetails/118-107th-Ave-NE-315-Bellevue-WA-98004/2070243355_zpid/',
    'https://www.zillow.com/b/9922-lake-washington-blvd-ne-bellevue-wa-9VSKqF/',
    'https://www.zillow.com/b/650-bellevue-way-ne-bellevue-wa-5XkDKg/',
    'https://www.zillow.com/b/washington-square-towers-bellevue-wa-5XkKmT/',
    'https://www.zillow.com/b/205-105th-ave-se-bellevue-wa-BzDfJd/',
    'https://www.zillow.com/homedetails/10650-NE-9th-Pl-UNIT-1524-Bellevue-WA-98004/2087878514_zpid/',
    'https://www.zillow.com/homedetails/1100-106th-Ave-NE-APT-606-Bellevue-WA-98004/58387391_zpid/',
    'https://www.zillow.com/homedetails/1620-103rd-Ave-NE-Bellevue-WA-98004/2104652731_zpid/',
    'https://www.zillow.com/b/188-bellevue-way-ne-bellevue-wa-BG9kPB/',
    'https://www.zillow.com/homedetails/1188-106th-Ave-NE-APT-321-Bellevue-WA-98004/60332740_zpid/',
    'https://www.zillow.com/homedetails/111-108th-Ave-NE-UNIT-A411-Bellevue-WA-98004/2070334836_zpid/',
    'https://www.zillow.com/b/812-100th-ave-ne-bellevue-wa-5jDpzz/',
    'https://www.zillow.com/homedetails/10700-NE-4th-St-UNIT-428-Bellevue-WA-98004/89210707_zpid/',
    'https://www.zillow.com/homedetails/10008-NE-16th-Pl-ADU-Bellevue-WA-98004/2072319855_zpid/',
    'https://www.zillow.com/homedetails/424-102nd-Ave-SE-APT-308-Bellevue-WA-98004/48966639_zpid/',
    'https://www.zillow.com/homedetails/550-100th-Ave-16-Bellevue-WA-98004/2070509346_zpid/',
    'https://www.zillow.com/b/10226-se-6th-street-bellevue-wa-5XwSLf/',
    'https://www.zillow.com/homedetails/10042-Main-St-APT-402-Bellevue-WA-98004/48957122_zpid/',
    'https://www.zillow.com/homedetails/10608-NE-2nd-St-FLOOR-12-ID117-Bellevue-WA-98004/2071029451_zpid/',
    'https://www.zillow.com/homedetails/10201-SE-3rd-St-1-Bellevue-WA-98004/2070980862_zpid/',
    'https://www.zillow.com/homedetails/425-Bellevue-Way-SE-25-Bellevue-WA-98004/2071064509_zpid/',
    'https://www.zillow.com/homedetails/425-Bellevue-Way-SE-62-Bellevue-WA-98004/2071064508_zpid/',
    'https://www.zillow.com/homedetails/417-99th-Ave-NE-G-Bellevue-WA-98004/2071688718_zpid/',
    'https://www.zillow.com/b/909-112th-ave-ne-bellevue-wa-9gyTGg/',
    'https://www.zillow.com/b/10608-ne-2nd-st-bellevue-wa-9q2SX2/',
    'https://www.zillow.com/homedetails/125-108th-Ave-SE-APT-12-Bellevue-WA-98004/2111358457_zpid/'
)

$indices = @(
    0,
    1,
    2,
    3,
    4,
    5,
    6,
    7,
    8,
    9,
    10,
    11,
    12,
    13,
    14,
    15,
    16,
    17,
    18,
    19,
    20,
    21,
    22,
    23,
    24,
    25,
    26,
    27,
    28,
    29,
    30,
    31,
    32,
    33,
    34,
    35,
    36,
    37,
    38,
    39,
    0,
    1,
    2,
    3,
    4,
    5,
    6,
    7,
    8,
    9,
    10,
    11,
    12,
    13,
    14,
    15,
    16,
    17,
    18,
    19,
    20,
    21
)

for ($n = 0; $n -lt $addresses.Length; $n++) {
    $row = $n + 2
    $ws.Cells.Item($row, 1).Value = $indices[$n]
    $ws.Cells.Item($row, 2).Value = $addresses[$n]
    $ws.Cells.Item($row, 3).Value = $prices[$n]
    $ws.Range("D$row").Hyperlinks.Delete()
    $ws.Range("D$row").Value = $urls[$n]
    $ws.Hyperlinks.Add($ws.Range("D$row"), $urls[$n])
    $ws.Range("D$row").Style = "Hyperlink"
}

